# "ebay add to cart" - add a second set of login credentials (email/password)
# to the "login_ebay" sheet, mirroring the existing first row (A2/B2), plus a
# third pair of rows styled like the small "Roboto" caption row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login_ebay")

# New validation email/password pair (row 3) -- styled like the existing
# hyperlinked email in A2.
$ws.Range("A3").Value = "anoor37290@gmail.com"
$ws.Range("B3").Value = "Serpent9999+"

# New email/password pair (row 4) -- styled like the small caption row B2.
$ws.Range("A4").Value = "justin.smith1@gmail.com"
$ws.Range("B4").Value = "AydenLiam1213@"

# Give A3 the same "Hyperlink" look as A2, then wire up a real mailto: link.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:anoor37290@gmail.com")

# Re-apply A2's format (PasteSpecial above can get clobbered by Hyperlinks.Add
# re-stamping its own font), so A3 keeps reusing the existing Hyperlink style.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# A4/B4 reuse the same small-caption formatting already used by B2.
$ws.Range("B2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

$excel.CutCopyMode = 0
